$d = $word.ActiveDocument

# --- Helpers -----------------------------------------------------------
# Append a brand-new empty paragraph right after $para and return it.
function AppendParaAfter($para) {
    $r = $para.Range
    $r.Collapse(1)          # wdCollapseEnd
    $r.InsertParagraphAfter()
    return $para.Next()
}

# Add plain text into a paragraph that is still completely empty.
function SetPlainText($para, [string]$text) {
    $para.Range.InsertAfter($text)
}

# Append a hyperlink run at the very end of $para (which may already
# contain a leading text run). Uses a one-character placeholder that is
# then replaced in-place by Hyperlinks.Add so that no stray empty <w:r/>
# is left behind afterwards. (Range.End here sits just past the
# paragraph mark, so the placeholder char is at End-2, not End-1.)
# The relationship Target is stored percent-encoded (matching the rest
# of this document's existing hyperlinks), while the visible run text
# keeps the literal, human-readable URL.
function AppendHyperlink($para, [string]$displayUrl) {
    $addr = $displayUrl.Replace("#", "%23")
    $para.Range.InsertAfter("X")
    $phPos = $para.Range.End - 2
    $hr = $d.Range($phPos, $phPos + 1)
    $d.Hyperlinks.Add($hr, $addr, $null, $null, $displayUrl) | Out-Null
}

# --- Locate the anchor paragraph ("...huawei-honor-5x.../#comment-8752") ---
$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*huawei-honor-5x-best-fingerprint-honor-holly2-plus-budget-phone*") {
        $anchor = $p
    }
}

$cur = $anchor

# 1) "4/1/2016" paragraph
$cur = AppendParaAfter $cur
SetPlainText $cur "4/1/2016"

# 2) "Awaiting approve " + hyperlink
$cur = AppendParaAfter $cur
SetPlainText $cur "Awaiting approve "
AppendHyperlink $cur "http://bloggingbehindthescenes.com/behind-the-scenes/guest-post-take-better-food-photos/#comment-1351"

# 3) hyperlink only
$cur = AppendParaAfter $cur
AppendHyperlink $cur "http://www.thesportsbank.net/bulls/jay-williams-believes-chicago-bulls-need-a-reboot-exclusive/#comment-3876929"

# 4) hyperlink only
$cur = AppendParaAfter $cur
AppendHyperlink $cur "http://allyouneedislists.com/design/banner-sizes/#comment-39489"

# 5) "Awaiting approve" (no trailing space)
$cur = AppendParaAfter $cur
SetPlainText $cur "Awaiting approve"

# 6) hyperlink only
$cur = AppendParaAfter $cur
AppendHyperlink $cur "http://migrationology.com/2015/11/lan-jia-gua-bao-taipei/#comment-291272"
